$wb = $excel.ActiveWorkbook

$wsFreq   = $wb.Worksheets.Item("FREQ")
$wsPontos = $wb.Worksheets.Item("PONTOS")

# New attendance values entered in column I (rows 3-20) of the FREQ sheet.
# "P" = presente, "F" = falta  (shared strings already used elsewhere in the sheet)
$iValues = [ordered]@{
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "F"
    8  = "F"
    9  = "P"
    10 = "F"
    11 = "P"
    12 = "F"
    13 = "P"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
}

foreach ($row in $iValues.Keys) {
    $wsFreq.Range("I$row").Value = $iValues[$row]
}

# PONTOS was the active/selected sheet before this edit; FREQ becomes the
# active sheet now, with I20 selected.
[void]$wsPontos.Range("D2").Select()
[void]$wsFreq.Activate()
[void]$wsFreq.Range("I20").Select()
